$d = $word.ActiveDocument

# Locate the paragraph containing the "aktivan - zastavica ..." bullet item so we can
# insert the new "oib - oib korisnika" bullet item right after it.
$targetIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*aktivan - zastavica*") {
        $targetIndex = $i
    }
    $i = $i + 1
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'aktivan - zastavica ...' paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)

# Insert a brand-new paragraph right after the target one. Word automatically clones
# the paragraph/list formatting (numPr numId=12, indentation, contextual spacing, rPr)
# from the paragraph it was split from, so the new bullet keeps the same list.
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "oib - oib korisnika"
